$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("I2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.3333333333333333
$ws.Range("Q2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.5
$ws.Range("U2").Value = 1
$ws.Range("W2").Value = 0.5
$ws.Range("C3").Value = 0.5
$ws.Range("F3").Value = 1
$ws.Range("I3").Value = 0.5
$ws.Range("M3").Value = 0.5
$ws.Range("Q3").Value = 0.5
$ws.Range("S3").Value = 0.5
$ws.Range("U3").Value = 0.5
$ws.Range("W3").Value = 1
$ws.Range("C4").Value = 0.6666666666666666
$ws.Range("F4").Value = 0.5
$ws.Range("I4").Value = 0.4
$ws.Range("M4").Value = 0.4
$ws.Range("Q4").Value = 0.2222222222222222
$ws.Range("S4").Value = 0.5
$ws.Range("U4").Value = 0.6666666666666666
$ws.Range("W4").Value = 0.6666666666666666
$ws.Range("C5").Value = 0.5555555555555556
$ws.Range("F5").Value = 0.7142857142857143
$ws.Range("I5").Value = 0.4545454545454545
$ws.Range("M5").Value = 0.4545454545454545
$ws.Range("Q5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.5
$ws.Range("U5").Value = 0.5555555555555556
$ws.Range("W5").Value = 0.8333333333333334
$ws.Range("C6").Value = 0.2754115523761866
$ws.Range("F6").Value = 0.6590018048024133
$ws.Range("I6").Value = 0.2754115523761866
$ws.Range("M6").Value = 0.2754115523761866
$ws.Range("Q6").Value = 0.2754115523761866
$ws.Range("S6").Value = 0.2754115523761866
$ws.Range("U6").Value = 0.2754115523761866
$ws.Range("W6").Value = 0.6885288809404666
